$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 2 (data row) string columns that were previously blank / placeholder
$ws.Range("C2").Value = "rt"
$ws.Range("D2").Value = "rt"
$ws.Range("E2").Value = "rt"
$ws.Range("F2").Value = "r"
$ws.Range("G2").Value = "tr"
$ws.Range("H2").Value = "t"
$ws.Range("I2").Value = "rt"
$ws.Range("J2").Value = "r"
$ws.Range("K2").Value = "t"
$ws.Range("L2").Value = "rt"
$ws.Range("M2").Value = "r"

# Numeric flag columns switched from 0 to 1
$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 1

$ws.Range("V2").Value = "rt"
$ws.Range("W2").Value = "rt"
